# Generate Report for Handoff
# The ecf43ef3-1d62-4099-80e8-e264b6095bb6.md file moved from
# "Handed back: in sync with en-US" to "Ready for handoff" status, with a
# refreshed handoff timestamp and a new "version not latest" error detail.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_dede = $wb.Worksheets.Item("de-de")

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a970326d18baa02626e29b77264306b1f302106f/e2e/ecf43ef3-1d62-4099-80e8-e264b6095bb6.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1ca43726ef90bde4f6f1b4db7c39a7aed0b627da/e2e/ecf43ef3-1d62-4099-80e8-e264b6095bb6.md."

# --- Overview sheet: row 3 is the ecf43ef3.md file ---
$ws_overview.Range("E3").Value = "Ready for handoff"
$ws_overview.Range("F3").Value = "Ready for handoff"
$ws_overview.Range("G3").Value = "2016-08-22 04:57:18"

# --- zh-cn sheet: row 3 is the ecf43ef3.md file ---
$ws_zhcn.Range("C3").Value = "Ready for handoff"
$ws_zhcn.Range("H3").Value = "2016-08-22 04:57:13"
$ws_zhcn.Range("P3").Value = $errorDetail
$ws_zhcn.Columns.Item(16).ColumnWidth = 39.16666666666667

# --- de-de sheet: row 3 is the ecf43ef3.md file ---
$ws_dede.Range("C3").Value = "Ready for handoff"
$ws_dede.Range("H3").Value = "2016-08-22 04:57:18"
$ws_dede.Range("P3").Value = $errorDetail
$ws_dede.Columns.Item(16).ColumnWidth = 39.16666666666667
